$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("#system")
$ws.Cells.Item(1,32).Value = "x"
$ws.Cells.Item(1,32).Value = ""
